# Add a new "2022-Q3" quarter to the 拼多多 holdings workbook.
#
# What this does, logically:
#   1. "总计" (summary) sheet gets a new row 2 for 2022-Q3 (1 fund held,
#      0.01 billion CNY), and the existing rows shift down one row with
#      their "持有数量(只)" sequence number (column A) incremented by one.
#   2. A brand-new detail sheet named "2022-Q3" is inserted right after
#      "总计" (i.e. before the existing "2022-Q2" sheet), holding a single
#      fund row with the new quarter's holding data.
#
# All of the other existing quarter tabs (2022-Q2, 2022-Q1, 2021-Q3,
# 2021-Q2, 2021-Q1) keep their own name and data untouched - they simply
# end up one tab position further to the right because of the newly
# inserted "2022-Q3" tab.

$wb = $excel.ActiveWorkbook

# Helper: write $val into $addr as literal TEXT (not auto-converted to a
# number) and strip the "quote prefix" style marker Excel normally adds,
# so the saved cell has no explicit style index (matches a plain data
# cell in this workbook).
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Update the "总计" summary sheet.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Insert a fresh row 2 above the current "2022-Q2" row; existing rows 2-6
# shift down to 3-7, formulas/data untouched.
$summary.Rows(2).Insert()

# The insert copied row-1 (header) formatting into the new row 2. Re-copy
# the correct per-column styles from row 3 (a genuine data row) instead.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B3:D3").Copy()
$summary.Range("B2:D2").PasteSpecial(-4122)

# Fill in the new 2022-Q3 summary row.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.01

# Renumber the sequence index in column A for the rows that shifted down
# (they kept their old 0..4 values; they need to become 1..5).
for ($r = 3; $r -le 7; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" detail sheet (same layout as "2022-Q2").
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

Set-TextValue $q3 "B2" "539002"
$q3.Range("C2").Value = "建信新兴市场优选混合（QDII）"
Set-TextValue $q3 "D2" "0.14"
Set-TextValue $q3 "E2" "81.57"
Set-TextValue $q3 "F2" "9.26"
Set-TextValue $q3 "G2" "0.0130"
$q3.Range("H2").Value = 2

# Restore the "总计" sheet as the active one (matches the original file).
$summary.Activate()
